$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 ("description") content updates -------------------------------

# New description cells G2 / H2 - clone the look of the existing
# "Neutral" (B2) / "Bad" (A2) header cells so the new cells share the same
# style entries, then fill in their text.
$ws2.Range("B2").Copy()
$ws2.Range("G2").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A2").Copy()
$ws2.Range("H2").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("G2").Value = "เป็นค่าว่างได้"
$ws2.Range("H2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

# New merged, centred header cell above them (G1:H1).
$ws2.Range("G1").HorizontalAlignment = -4108   # xlHAlignCenter
$ws2.Range("H1").HorizontalAlignment = -4108   # xlHAlignCenter
$ws2.Range("G1:H1").Merge()
$ws2.Columns("G:H").ColumnWidth = 24.8

# Update the existing description text for taxrate_id (B1) to include the
# extra "numbers only" note.
$ws2.Range("B1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# --- Sheet1 selection (no longer the active tab) ---------------------------
$ws1.Activate()
$ws1.Range("A11").Select()

# --- Sheet2 becomes the active tab with its own selection -------------------
$ws2.Activate()
$ws2.Range("F11").Select()
